# (#33) Alteração nos rótulos da tabela para já transformar a primeira linha
# em cabeçalho automaticamente no Power BI.
# Prefix the year/interval header labels in row 1 of each sheet so that
# Power BI can auto-detect the first row as a header.

$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3, 5: single-year headers -> prefix "Ano "
$anoSheets = @(
  "Potencia Acumulada - SIN (MW)",
  "Geracao Periodo Medio (MWMed)",
  "Atendimento a Ponta(MW)",
  "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
  $ws = $wb.Worksheets.Item($sheetName)
  $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Value2
  $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Value2
  $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Value2
  $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Value2
}

# Sheet 4: interval headers -> prefix "Intervalo "
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$ws4.Range("B1").Value = "Intervalo " + $ws4.Range("B1").Value2
$ws4.Range("C1").Value = "Intervalo " + $ws4.Range("C1").Value2
$ws4.Range("D1").Value = "Intervalo " + $ws4.Range("D1").Value2
$ws4.Range("E1").Value = "Intervalo " + $ws4.Range("E1").Value2

# Sheet 6: only a single year column (B1) -> prefix "Ano "
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Range("B1").Value = "Ano " + $ws6.Range("B1").Value2
